$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 926487.6
$ws.Range("J17").Value = 1089632.6
$ws.Range("L17").Value = 3268897.8
$ws.Range("N17").Value = -3269233.8

$ws.Range("H99").Value = 1011.5294
$ws.Range("I99").Value = 313
$ws.Range("J99").Value = 2688
$ws.Range("K99").Value = 939
$ws.Range("L99").Value = 8064
$ws.Range("M99").Value = 559
$ws.Range("N99").Value = -11060

$ws.Range("H116").Value = 10950.808
$ws.Range("I116").Value = 7804.9165
$ws.Range("J116").Value = 13647.286
$ws.Range("K116").Value = 7804.9165
$ws.Range("L116").Value = 13647.286
$ws.Range("M116").Value = -4362.9165
$ws.Range("N116").Value = -20531.286

$ws.Range("H131").Value = 6105.7646
$ws.Range("I131").Value = 1986.625
$ws.Range("K131").Value = 5959.875
$ws.Range("M131").Value = -919.875

$ws.Range("H135").Value = 17546298
$ws.Range("I135").Value = 25643126
$ws.Range("J135").Value = 3170.6667
$ws.Range("K135").Value = 230788134
$ws.Range("L135").Value = 28536.0003
$ws.Range("M135").Value = -230785599
$ws.Range("N135").Value = -33606.0003

$ws.Range("H137").Value = 722767.9
$ws.Range("I137").Value = 17977.875
$ws.Range("K137").Value = 53933.625
$ws.Range("M137").Value = -51383.625

$ws.Range("H138").Value = 4220.375
$ws.Range("I138").Value = 1943.7
$ws.Range("J138").Value = 4587.5806
$ws.Range("K138").Value = 5831.1
$ws.Range("L138").Value = 13762.7418
$ws.Range("M138").Value = -691.1000000000004
$ws.Range("N138").Value = -24042.7418

$ws.Range("H141").Value = 6900
$ws.Range("I141").Value = 6800
$ws.Range("K141").Value = 20400
$ws.Range("M141").Value = -15220

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2040.3077
$ws.Range("I2").Value = 2390.5
$ws.Range("J2").Value = 1480
$ws.Range("K2").Value = 2390.5
$ws.Range("L2").Value = 1480
$ws.Range("M2").Value = -2277.5
$ws.Range("N2").Value = -1706

$ws.Range("H32").Value = 3383.5476
$ws.Range("I32").Value = 3080.743
$ws.Range("K32").Value = 3080.743
$ws.Range("M32").Value = -2793.743

$ws.Range("H116").Value = 2040.3077
$ws.Range("I116").Value = 2390.5
$ws.Range("J116").Value = 1480
$ws.Range("K116").Value = 2390.5
$ws.Range("L116").Value = 1480
$ws.Range("M116").Value = -96.5
$ws.Range("N116").Value = -6068

$ws.Range("H122").Value = 4697.5
$ws.Range("I122").Value = 5333.3335
$ws.Range("K122").Value = 16000.0005
$ws.Range("M122").Value = -13550.0005

$ws.Range("H132").Value = 3024.913
$ws.Range("I132").Value = 2578.65
$ws.Range("K132").Value = 7735.950000000001
$ws.Range("M132").Value = -5205.950000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2040.3077
$ws.Range("I3").Value = 2390.5
$ws.Range("J3").Value = 1480
$ws.Range("K3").Value = 2390.5
$ws.Range("L3").Value = 1480
$ws.Range("M3").Value = -2276.5
$ws.Range("N3").Value = -1708

$ws.Range("H94").Value = 1544.5625
$ws.Range("I94").Value = 1580.6154
$ws.Range("K94").Value = 1580.6154
$ws.Range("M94").Value = -1129.6154

$ws.Range("H140").Value = 71643.89999999999
$ws.Range("J140").Value = 71643.89999999999
$ws.Range("L140").Value = 71643.89999999999
$ws.Range("N140").Value = -82003.89999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2485470.5
$ws.Range("I31").Value = 3089596.5
$ws.Range("K31").Value = 3089596.5
$ws.Range("M31").Value = -3089301.5

$ws.Range("H34").Value = 2485470.5
$ws.Range("I34").Value = 3089596.5
$ws.Range("K34").Value = 3089596.5
$ws.Range("M34").Value = -3089394.5

$ws.Range("H107").Value = 1368.1538
$ws.Range("I107").Value = 1453.4166
$ws.Range("K107").Value = 1453.4166
$ws.Range("M107").Value = 466.5834

$ws.Range("H122").Value = 1884.0834
$ws.Range("I122").Value = 2296.125
$ws.Range("J122").Value = 1060
$ws.Range("K122").Value = 6888.375
$ws.Range("L122").Value = 3180
$ws.Range("M122").Value = -4438.375
$ws.Range("N122").Value = -8080

$ws.Range("H132").Value = 1693.0513
$ws.Range("I132").Value = 1553.3948
$ws.Range("K132").Value = 4660.1844
$ws.Range("M132").Value = -2130.1844

$ws.Range("H134").Value = 19421.486
$ws.Range("I134").Value = 21282.484
$ws.Range("K134").Value = 63847.452
$ws.Range("M134").Value = -61312.452

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 25641918
$ws.Range("J107").Value = 2228.3333
$ws.Range("L107").Value = 6684.999899999999
$ws.Range("N107").Value = -10524.9999

$ws.Range("H129").Value = 1325.375
$ws.Range("I129").Value = 681.8
$ws.Range("J129").Value = 2398
$ws.Range("K129").Value = 2045.4
$ws.Range("L129").Value = 7194
$ws.Range("M129").Value = 2954.6
$ws.Range("N129").Value = -17194

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 168
$ws.Range("I2").Value = 132.18182
$ws.Range("J2").Value = 266.5
$ws.Range("K2").Value = 132.18182
$ws.Range("L2").Value = 266.5
$ws.Range("M2").Value = -19.18181999999999
$ws.Range("N2").Value = -492.5

$ws.Range("H97").Value = 4891
$ws.Range("I97").Value = 5158.5713
$ws.Range("J97").Value = 4266.6665
$ws.Range("K97").Value = 5158.5713
$ws.Range("L97").Value = 4266.6665
$ws.Range("M97").Value = -4662.5713
$ws.Range("N97").Value = -5258.6665

$ws.Range("H102").Value = 7499.5
$ws.Range("I102").Value = 7499.5
$ws.Range("K102").Value = 7499.5
$ws.Range("M102").Value = -5877.5

$ws.Range("H122").Value = 8120.1113
$ws.Range("I122").Value = 8120.1113
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 24360.3339
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -21910.3339
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 16635.63
$ws.Range("I132").Value = 18186.088
$ws.Range("K132").Value = 54558.264
$ws.Range("M132").Value = -52028.264

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11349.846
$ws.Range("I7").Value = 13634.223
$ws.Range("K7").Value = 13634.223
$ws.Range("M7").Value = -13522.223

$ws.Range("H126").Value = 11349.846
$ws.Range("I126").Value = 13634.223
$ws.Range("K126").Value = 40902.669
$ws.Range("M126").Value = -38432.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5540.886
$ws.Range("I132").Value = 4366.143
$ws.Range("K132").Value = 13098.429
$ws.Range("M132").Value = -10568.429

$ws.Range("H136").Value = 1360.8214
$ws.Range("I136").Value = 1084.16
$ws.Range("K136").Value = 3252.48
$ws.Range("M136").Value = -702.4800000000005
